# Sort the calibration data (rows 2-12) ascending by column A (time),
# keeping the header row (row 1) fixed in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D12")
$keyRange  = $ws.Range("A2:A12")

$dataRange.Sort($keyRange, 1)
